# Applies the "Updates to model policy subscripts" commit to the
# BAEPAbCiPC workbook: refreshes the About-tab narrative text and
# relabels a couple of "not used" fuel rows on the lever tab.

$wb  = $excel.ActiveWorkbook
$about = $wb.Worksheets.Item(1)      # "About" sheet
$lever = $wb.Worksheets.Item(2)      # "BAEPAbCiPC" sheet

# ---------------------------------------------------------------
# About sheet: rewrite the explanatory paragraph (rows 11-13) and
# append a brand-new second paragraph (rows 15-18). Row 14 is left
# blank, same as the gap that separates rows 9/11 and 13/15.
# ---------------------------------------------------------------
$about.Range("A11").Value = "As of EPS 3.1.0, this lever supports the three energy carriers (electricity,"
$about.Range("A12").Value = "district heat, and hydrogen), as well as fuels produced by the natural gas"
$about.Range("A13").Value = "and petroleum, coal, biomass, and biofuel industries, as noted on the blue tab."
$about.Range("A14").Value = ""

# ---------------------------------------------------------------
# Lever sheet: relabel the two rows that are kept for documentation
# only but are no longer active subscript members.
# ---------------------------------------------------------------
$lever.Range("A6:B6").Copy() | Out-Null
$lever.Range("A21:B21").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$lever.Range("A21").Value = "municipal solid waste (NOT USED)"

$lever.Range("A6:B6").Copy() | Out-Null
$lever.Range("A5:B5").PasteSpecial(-4122) | Out-Null      # xlPasteFormats
$lever.Range("A5").Value = "nuclear (NOT USED)"

# Header cell: "Boolean" -> "Unit: boolean (1 or 0)", now styled in
# italics instead of the old header style.
$lever.Range("A1").Value = "Unit: boolean (1 or 0)"
$lever.Range("A1").Font.Italic = $true

$excel.CutCopyMode = 0

# ---------------------------------------------------------------
# New paragraph 2 on the About sheet (rows 15-18), added after the
# lever-sheet edits so the shared-string table append order matches
# the authored workbook.
# ---------------------------------------------------------------
$about.Range("A15").Value = "In the U.S. model, by default, we allow the suppliers of energy carriers"
$about.Range("A16").Value = "(electricity, district heat, and hydrogen) to pass through changes in their"
$about.Range("A17").Value = "expenses, while other fuel suppliers do not, due to the influence of a global"
$about.Range("A18").Value = "market on setting prices."
